# Temperature controller for DH thermal storage
$wb = $excel.ActiveWorkbook

$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsStorage    = $wb.Worksheets.Item("storageAssets")

# --- conversionAssets (sheet3) ---
# F8: 1000 -> 300
$wsConversion.Range("F8").Value = 300

# --- storageAssets (sheet4) ---
# Replace formulas in F3:F9 (=20*100/1000, =30*100/1000, ...) with a flat value of 100
$wsStorage.Range("F3:F9").Value = 100

# F11: 100 -> 1000
$wsStorage.Range("F11").Value = 1000

# --- sheet view / selection updates ---
# storageAssets ends up with selection F11 (not the active tab)
$wsStorage.Activate()
$wsStorage.Range("F11").Select()

# conversionAssets becomes the active tab with selection F8
$wsConversion.Activate()
$wsConversion.Range("F8").Select()
